# StateFunction.xlsx: "unify the conception of DataNode, DataTable, Entity."
# The sheet that used to represent a generic "Property" table is renamed to
# "DataNode", and the active selection is left on D26 (matching the editor's
# on-disk cursor position after the rename).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "Property1" to "DataNode".
$ws.Name = "DataNode"

# Make sure we're looking at the renamed sheet and park the selection on D26,
# same as the saved workbook state.
$ws.Activate()
$ws.Range("D26").Select()
